$d = $word.ActiveDocument

# The document's String collection used to keep empty-text paragraphs
# (just a blank run) mixed in with the real lines. Drop those blank
# paragraphs so only the populated strings remain.
# Walk back-to-front so deleting a paragraph never shifts the index of
# one we haven't visited yet.
$count = $d.Paragraphs.Count
for ($idx = $count; $idx -ge 1; $idx--) {
    $p = $d.Paragraphs($idx)
    $text = $p.Range.Text
    # A paragraph's Range.Text includes the trailing paragraph mark
    # (length 1) when there is no other content.
    if ($text.Length -le 1) {
        $p.Range.Delete()
    }
}
